# "changes to power measurements"
# Adds a new "RECEIVE DATA_OVERHEAR" block (columns I/J), extends the
# V_BATTERY measurement rows (5-8) with #3/#4/#5 Avg. data for 2.8V and
# 2.7V battery levels, recolors the relevant blocks, and repositions the
# selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New / changed cell VALUES
# ---------------------------------------------------------------------

# New "RECEIVE DATA_OVERHEAR" header columns
$ws.Range("I1").Value = "RECEIVE DATA_OVERHEAR_V"
$ws.Range("J1").Value = "RECEIVE DATA_OVERHEAR_TIME"

# #1 Avg. overhear figures
$ws.Range("I3").Value = "196 mV"
$ws.Range("J3").Value = "3,6 ms"

# #2 Avg. overhear figures
$ws.Range("I4").Value = "211 mV"
$ws.Range("J4").Value = "3 ms"

# #3 Avg. row (2,8 V_BATTERY block)
$ws.Range("B5").Value = "187 mV"
$ws.Range("C5").Value = "9,4 ms"

# #4 Avg. row (start of 2,8 V_BATTERY label + its figures)
$ws.Range("B6").Value = "2,8 V_BATTERY"
$ws.Range("D6").Value = "202 mV"
$ws.Range("E6").Value = "216 mV"
$ws.Range("F6").Value = "9 ms"

# #5 Avg. row
$ws.Range("D7").Value = "205 mV"
$ws.Range("E7").Value = "216 mV"
$ws.Range("F7").Value = "10,5 ms"

# 2,7 V_BATTERY figure
$ws.Range("F8").Value = "2,7 V_BATTERY"

# Note row
$ws.Range("E17").Value = "Probably only overhear"

# ---------------------------------------------------------------------
# 2. Fills / styles
# ---------------------------------------------------------------------

# Fill "2" (theme 7 / accent4, tinted) - lightly highlighted existing
# V_SHUNT block + the new note cell
# (applied per contiguous block - multi-area Range assignment only
# reliably touches the first area on this engine)
$ws.Range("B3:F4").Interior.ThemeColor = 8
$ws.Range("D5").Interior.ThemeColor = 8
$ws.Range("F5").Interior.ThemeColor = 8
$ws.Range("E17").Interior.ThemeColor = 8

# Fill "3" (theme 9 / accent6) - the new OVERHEAR columns + the new
# V_BATTERY rows
$ws.Range("I2:J4").Interior.ThemeColor = 10
$ws.Range("B5:C5").Interior.ThemeColor = 10
$ws.Range("B6:F6").Interior.ThemeColor = 10
$ws.Range("D7:F7").Interior.ThemeColor = 10
$ws.Range("D8:F8").Interior.ThemeColor = 10
$ws.Range("B11").Interior.ThemeColor = 10
$ws.Range("B13:B14").Interior.ThemeColor = 10

# Header cells for the new OVERHEAR columns are bold (fontId 1) as well
# as filled
$ws.Range("I1:J1").Font.Bold = $true
$ws.Range("I1:J1").Interior.ThemeColor = 10

# B12 keeps its 2-decimal number format but also gets fill "2"
$ws.Range("B12").Interior.ThemeColor = 8
$ws.Range("B12").NumberFormat = "0.00"

# ---------------------------------------------------------------------
# 3. Column widths
# ---------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 19.166666666666668
$ws.Columns.Item(9).ColumnWidth = 24.666666666666668
$ws.Columns.Item(10).ColumnWidth = 27.666666666666668

# ---------------------------------------------------------------------
# 4. Selection
# ---------------------------------------------------------------------
$ws.Range("D13").Select()
